# Adding 29/12/2017 workout data
#
# Appends 6 new workout-log rows (exercise date 29/12/2017, a Friday in
# December 2017, week 52) to the bottom of the WeightTraining sheet:
#   Barbell Squat, Deadlift, Leg Extension, Leg Curl, Hip adduction,
#   Hip abduction.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=ExerciseId B=DateId C=Exercise Date D=Exercise Week
#          E=Exercise Month F=Exercise Year G=Exercise Day H=Exercise Name
#          I=Weight J=Sets K=Reps
$newRows = @(
    @(307, 36, 43098, 52, "December", 2017, "Friday", "Barbell Squat",  90,  5,  5),
    @(308, 36, 43098, 52, "December", 2017, "Friday", "Deadlift",      130,  5,  5),
    @(309, 36, 43098, 52, "December", 2017, "Friday", "Leg Extension", 100,  4, 12),
    @(310, 36, 43098, 52, "December", 2017, "Friday", "Leg Curl",       60,  4, 12),
    @(311, 36, 43098, 52, "December", 2017, "Friday", "Hip adduction",  54,  4, 12),
    @(312, 36, 43098, 52, "December", 2017, "Friday", "Hip abduction",  45,  4, 12)
)

$firstNewSheetRow = 308
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $sheetRow = $firstNewSheetRow + $i
    $values = $newRows[$i]
    for ($col = 0; $col -lt $values.Count; $col++) {
        $ws.Cells.Item($sheetRow, $col + 1).Value = $values[$col]
    }
}

# Move the active selection to just below the newly added data, as in the
# edited workbook.
[void]$ws.Range("A315").Select()
